# Remove start/end datetimes from IMC metadata template. Fix #561.
#
# Column "AF" (end_datetime) and column "AM" (start_datetime) are removed
# entirely from the "Export as TSV" sheet. Removing them shifts every
# column after each deletion point one slot to the left (comments,
# data validations, and the header row included).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# Delete the two columns, right-to-left, so earlier column letters stay
# valid while later ones are being removed.
$ws.Range("AM1").EntireColumn.Delete()
$ws.Range("AF1").EntireColumn.Delete()

# The engine's EntireColumn.Delete() shifts cell values/styles/validations
# but leaves cell comments anchored to their original addresses, so fix the
# comments up by hand to match what real Excel would have produced: every
# comment from the old AG1..AP1 range slides left by two columns into the
# new AF1..AN1 range, and the two trailing (now out-of-range) comments are
# removed.
$null = $ws.Range("AF1").Comment.Text("Image width value of the ROI acquisition")
$null = $ws.Range("AG1").Comment.Text("Units of image width of the ROI acquisition")
$null = $ws.Range("AH1").Comment.Text("Image height value of the ROI acquisition")
$null = $ws.Range("AI1").Comment.Text("Units of image height of the ROI acquisition")
$null = $ws.Range("AJ1").Comment.Text('This refers to the data type, which is a "float" for the IMC counts.')
$null = $ws.Range("AK1").Comment.Text("Type of signal measured per channel (usually dual counts)")
$null = $ws.Range("AL1").Comment.Text("Numerical data precision in bytes")
$null = $ws.Range("AM1").Comment.Text("Relative path to file with ORCID IDs for contributors for this dataset.")
$null = $ws.Range("AN1").Comment.Text("Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.")

# These two addresses are now past the used range (A1:AN1) -- drop their
# stale comments so nothing is left dangling off the end of the sheet.
$ws.Range("AO1").Comment.Delete()
$ws.Range("AP1").Comment.Delete()
